$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark the rows corresponding to the newly-implemented features as
# achieved on Milestone I ("I" in column E, "X" in column F):
#   row 6  - Geometry Instancing w/ 1 drawInstance related call
#   row 18 - Applying applicable color map texturing to drawn geometry
#   row 30 - Applying functional directional light to drawn geometry
#   row 31 - Applying functional point light to drawn geometry
#   row 33 - Combining 2 functional lights on the same drawn geometry
#   row 34 - Demostrates dynamic change in direction of directional lighting
#   row 35 - Demostrates dynamic change in position of point lighting
#   row 56 - Simultaneous Loading of Textures or Models done with Multithreading
$rows = @(6, 18, 30, 31, 33, 34, 35, 56)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "I"
    $ws.Range("F$r").Value = "X"
}

# Leave a note for the grader/reviewer on row 73, matching the style of
# the existing "Author:" note already present on row 68.
$comment = $ws.Range("A73").AddComment("Author:" + [char]10)
$null = $comment

# Update the view: scroll position reset and a new active selection cell.
$null = $ws.Range("E26").Select()
